# TourPlanner_Checklist.xlsx - "update todos - checklist einschätzung"
#
# Adds a new column E ("unsere Einschätzung" = "our assessment") next to the
# existing Max. Points column, fills in point estimates for most of the
# checklist items, and sums them in E78 (mirrors the existing C78 SUM).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header for the new column ------------------------------------------------
$ws.Range("E26").Value = "unsere Einschätzung"

# --- "Must Haves" section values ---------------------------------------------
$ws.Range("E27").Value = 3.5
$ws.Range("E28").Value = 0.5
$ws.Range("E29").Value = 1.5

# --- "Tours" section values ---------------------------------------------------
$ws.Range("E32").Value = 2
$ws.Range("E33").Value = 1.5
$ws.Range("E34").Value = 1.5
$ws.Range("E35").Value = 2

# These four rows don't already carry a bold row style, so the new cells need
# the bold/size-12 font applied explicitly (matches the header-row style used
# elsewhere in the sheet, e.g. A26/A31/A38).
$boldRange = $ws.Range("E32:E35")
$boldRange.Font.Bold = $true
$boldRange.Font.Size = 12

# --- "Tour Logs" section values -----------------------------------------------
$ws.Range("E39").Value = 3
$ws.Range("E40").Value = 2
$ws.Range("E41").Value = 2
$ws.Range("E42").Value = 1

# --- "Reports, Import/Export" section values ----------------------------------
$ws.Range("E49").Value = 2
$ws.Range("E50").Value = 1
$ws.Range("E51").Value = 1.5
$ws.Range("E52").Value = 1.5

# --- "Mandatory Unique feature" -----------------------------------------------
$ws.Range("E54").Value = 2

# --- "Non-Functional Requirements" section values -----------------------------
$ws.Range("E57").Value = 4
$ws.Range("E59").Value = 2
$ws.Range("E60").Value = 2
$ws.Range("E61").Value = 0.5

# --- Total row: mirror the existing C78 SUM formula for the new column -------
$ws.Range("E78").Formula = "=SUM(E27:E74)"

# --- Leave the selection where the author ended up editing --------------------
$ws.Range("E27").Select()
